# Data refresh: update market price / profit figures across crafting-job sheets
# (values sourced from the latest Universalis market snapshot).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: H51, I51, J51, K51, L51, M51, N51
$ws.Range("H51").Value = 6333.3335
$ws.Range("I51").Value = 6000
$ws.Range("J51").Value = 7000
$ws.Range("K51").Value = 6000
$ws.Range("L51").Value = 7000
$ws.Range("M51").Value = -5516
$ws.Range("N51").Value = -7968
# Row 62: H62, I62, K62, M62
$ws.Range("H62").Value = 166673490
$ws.Range("I62").Value = 200002580
$ws.Range("K62").Value = 200002580
$ws.Range("M62").Value = -200001956
# Row 65: H65, I65, K65, M65
$ws.Range("H65").Value = 166673490
$ws.Range("I65").Value = 200002580
$ws.Range("K65").Value = 1000012900
$ws.Range("M65").Value = -1000009780
# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value = 1145
$ws.Range("I107").Value = 627.4167
$ws.Range("J107").Value = 2387.2
$ws.Range("K107").Value = 627.4167
$ws.Range("L107").Value = 2387.2
$ws.Range("M107").Value = 1292.5833
$ws.Range("N107").Value = -6227.2
# Row 112: H112, J112, L112, N112
$ws.Range("H112").Value = 1402.9656
$ws.Range("J112").Value = 1402.9656
$ws.Range("L112").Value = 4208.8968
$ws.Range("N112").Value = -6424.8968
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 833.5161000000001
$ws.Range("I132").Value = 730.62964
$ws.Range("K132").Value = 2191.88892
$ws.Range("M132").Value = 338.1110800000001
# Row 137: H137, I137, J137, K137, L137, M137, N137
$ws.Range("H137").Value = 2265.8235
$ws.Range("I137").Value = 2060.5454
$ws.Range("J137").Value = 2642.1667
$ws.Range("K137").Value = 6181.6362
$ws.Range("L137").Value = 7926.500100000001
$ws.Range("M137").Value = -3631.6362
$ws.Range("N137").Value = -13026.5001
# Row 138: H138, I138, J138, K138, L138, M138, N138
$ws.Range("H138").Value = 1839.8983
$ws.Range("I138").Value = 1383.2174
$ws.Range("J138").Value = 2131.6667
$ws.Range("K138").Value = 4149.6522
$ws.Range("L138").Value = 6395.000100000001
$ws.Range("M138").Value = 990.3477999999996
$ws.Range("N138").Value = -16675.0001
# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 72284
$ws.Range("J140").Value = 72284
$ws.Range("L140").Value = 72284
$ws.Range("N140").Value = -82644

$ws = $wb.Worksheets.Item("ARM")
# Row 16: H16, J16, L16, N16
$ws.Range("H16").Value = 1899.5
$ws.Range("J16").Value = 1899.5
$ws.Range("L16").Value = 1899.5
$ws.Range("N16").Value = -2473.5
# Row 32: H32, I32, K32, M32
$ws.Range("H32").Value = 5418.8965
$ws.Range("I32").Value = 3408.7551
$ws.Range("K32").Value = 3408.7551
$ws.Range("M32").Value = -3121.7551
# Row 45: H45, I45, K45, M45
$ws.Range("H45").Value = 5295501
$ws.Range("I45").Value = 12857990
$ws.Range("K45").Value = 12857990
$ws.Range("M45").Value = -12857613
# Row 74: H74, I74, K74, M74
$ws.Range("H74").Value = 2867.5557
$ws.Range("I74").Value = 1502.75
$ws.Range("K74").Value = 1502.75
$ws.Range("M74").Value = -628.75
# Row 77: H77, I77, K77, M77
$ws.Range("H77").Value = 2867.5557
$ws.Range("I77").Value = 1502.75
$ws.Range("K77").Value = 7513.75
$ws.Range("M77").Value = -3145.75
# Row 110: H110, I110, K110, M110
$ws.Range("H110").Value = 1403.2222
$ws.Range("I110").Value = 275.8
$ws.Range("K110").Value = 275.8
$ws.Range("M110").Value = 1769.2
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 1192.6842
$ws.Range("I122").Value = 870.9167
$ws.Range("K122").Value = 2612.7501
$ws.Range("M122").Value = -162.7501000000002

$ws = $wb.Worksheets.Item("BSM")
# Row 94: H94, I94, J94, K94, L94, M94, N94
$ws.Range("H94").Value = 491.8
$ws.Range("I94").Value = 384
$ws.Range("J94").Value = 743.3333
$ws.Range("K94").Value = 384
$ws.Range("L94").Value = 743.3333
$ws.Range("M94").Value = 67
$ws.Range("N94").Value = -1645.3333
# Row 132: H132, J132, L132, N132
$ws.Range("H132").Value = 33500
$ws.Range("J132").Value = 33500
$ws.Range("L132").Value = 33500
$ws.Range("N132").Value = -43620
# Row 134: H134, I134, K134, M134
$ws.Range("H134").Value = 5532.4814
$ws.Range("I134").Value = 6066.091
$ws.Range("K134").Value = 18198.273
$ws.Range("M134").Value = -15663.273

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31, I31, J31, K31, L31, M31, N31
$ws.Range("H31").Value = 2364.7058
$ws.Range("I31").Value = 2433.3333
$ws.Range("J31").Value = 2350
$ws.Range("K31").Value = 2433.3333
$ws.Range("L31").Value = 2350
$ws.Range("M31").Value = -2138.3333
$ws.Range("N31").Value = -2940
# Row 34: H34, I34, J34, K34, L34, M34, N34
$ws.Range("H34").Value = 2364.7058
$ws.Range("I34").Value = 2433.3333
$ws.Range("J34").Value = 2350
$ws.Range("K34").Value = 2433.3333
$ws.Range("L34").Value = 2350
$ws.Range("M34").Value = -2231.3333
$ws.Range("N34").Value = -2754
# Row 58: H58, I58, K58, M58
$ws.Range("H58").Value = 5437582.5
$ws.Range("I58").Value = 14493486
$ws.Range("K58").Value = 14493486
$ws.Range("M58").Value = -14493283
# Row 99: H99, I99, K99, M99
$ws.Range("H99").Value = 1430633.6
$ws.Range("I99").Value = 2001787
$ws.Range("K99").Value = 2001787
$ws.Range("M99").Value = -2000289
# Row 126: H126, I126, K126, M126
$ws.Range("H126").Value = 1430633.6
$ws.Range("I126").Value = 2001787
$ws.Range("K126").Value = 6005361
$ws.Range("M126").Value = -6002891
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 2324.16
$ws.Range("I132").Value = 1449.7778
$ws.Range("K132").Value = 4349.3334
$ws.Range("M132").Value = -1819.3334
# Row 134: H134, I134, K134, M134
$ws.Range("H134").Value = 3102.3572
$ws.Range("I134").Value = 2780.7
$ws.Range("K134").Value = 8342.099999999999
$ws.Range("M134").Value = -5807.099999999999
# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 5437582.5
$ws.Range("I136").Value = 14493486
$ws.Range("K136").Value = 43480458
$ws.Range("M136").Value = -43477908

$ws = $wb.Worksheets.Item("CUL")
# Row 5: H5, J5, L5, N5
$ws.Range("H5").Value = 400.61905
$ws.Range("J5").Value = 599.2
$ws.Range("L5").Value = 1797.6
$ws.Range("N5").Value = -2021.6
# Row 22: H22, I22, J22, K22, L22, M22, N22
$ws.Range("H22").Value = 3214.6924
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 3232.5833
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 9697.749899999999
$ws.Range("M22").Value = -8831
$ws.Range("N22").Value = -10035.7499
# Row 27: H27, I27, J27, K27, L27, M27, N27
$ws.Range("H27").Value = 3214.6924
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 3232.5833
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 9697.749899999999
$ws.Range("M27").Value = -8898
$ws.Range("N27").Value = -9901.749899999999
# Row 131: H131, J131, L131, N131
$ws.Range("H131").Value = 11406.04
$ws.Range("J131").Value = 11697.575
$ws.Range("L131").Value = 35092.72500000001
$ws.Range("N131").Value = -45172.72500000001
# Row 135: H135, J135, L135, N135
$ws.Range("H135").Value = 400.61905
$ws.Range("J135").Value = 599.2
$ws.Range("L135").Value = 5392.8
$ws.Range("N135").Value = -10462.8
# Row 137: H137, I137, J137, K137, L137, M137, N137
$ws.Range("H137").Value = 3423.2
$ws.Range("I137").Value = 1231.4546
$ws.Range("J137").Value = 5145.2856
$ws.Range("K137").Value = 3694.3638
$ws.Range("L137").Value = 15435.8568
$ws.Range("M137").Value = 1405.6362
$ws.Range("N137").Value = -25635.8568

$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70
$ws.Range("H70").Value = 3919.8
# Row 73: H73
$ws.Range("H73").Value = 3919.8
# Row 102: H102, I102, J102, K102, L102, M102, N102
$ws.Range("H102").Value = 2775.8823
$ws.Range("I102").Value = 2561.3845
$ws.Range("J102").Value = 3473
$ws.Range("K102").Value = 2561.3845
$ws.Range("L102").Value = 3473
$ws.Range("M102").Value = -939.3845000000001
$ws.Range("N102").Value = -6717
# Row 126: H126, I126, K126, M126
$ws.Range("H126").Value = 2830524.2
$ws.Range("I126").Value = 5053509
$ws.Range("K126").Value = 15160527
$ws.Range("M126").Value = -15158057

$ws = $wb.Worksheets.Item("LTW")
# Row 61: H61, I61, J61, K61, L61, M61, N61
$ws.Range("H61").Value = 2477.0527
$ws.Range("I61").Value = 2355.3333
$ws.Range("J61").Value = 2685.7144
$ws.Range("K61").Value = 2355.3333
$ws.Range("L61").Value = 2685.7144
$ws.Range("M61").Value = -2153.3333
$ws.Range("N61").Value = -3089.7144
# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 2477.0527
$ws.Range("I113").Value = 2355.3333
$ws.Range("J113").Value = 2685.7144
$ws.Range("K113").Value = 2355.3333
$ws.Range("L113").Value = 2685.7144
$ws.Range("M113").Value = -185.3332999999998
$ws.Range("N113").Value = -7025.7144
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 1867.9773
$ws.Range("I132").Value = 1287.2632
$ws.Range("K132").Value = 3861.7896
$ws.Range("M132").Value = -1331.7896
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 4559.4
$ws.Range("I136").Value = 3264.1667
$ws.Range("J136").Value = 6502.25
$ws.Range("K136").Value = 9792.500100000001
$ws.Range("L136").Value = 19506.75
$ws.Range("M136").Value = -7242.500100000001
$ws.Range("N136").Value = -24606.75

$ws = $wb.Worksheets.Item("WVR")
# Row 126: H126, I126, K126, M126
$ws.Range("H126").Value = 1997.5
$ws.Range("I126").Value = 2009.5
$ws.Range("K126").Value = 6028.5
$ws.Range("M126").Value = -3558.5
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1091.0968
$ws.Range("I132").Value = 890.42224
$ws.Range("J132").Value = 1622.2941
$ws.Range("K132").Value = 2671.26672
$ws.Range("L132").Value = 4866.8823
$ws.Range("M132").Value = -141.2667200000001
$ws.Range("N132").Value = -9926.882300000001
